# The workbook's sheet1 cell A1 holds a status message built by a
# "check for past dates" report. In the new run, dates in the PAST were
# found, so the leading "No dates found in the PAST. " sentence needs to
# be dropped, leaving only the three "date ... is in the PAST" sentences.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "The date from 05/12/2022 is in the PAST. The date from 08/12/2022 is in the PAST. The date from 21/11/2022 is in the PAST"
